$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: JTSJOL_mil - Job Openings - updated to Nov 2025 data
$ws.Range("E7").Value = 7.146
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "Nov 2025"
$ws.Range("G7").Value = 7.791380165289254
$ws.Range("H7").Value = -0.8850000000000007
$ws.Range("I7").Value = -0.1101979828165858

# Row 9: ICSA_thou - Initial Jobless Claims - updated to Jan 2026 data
$ws.Range("E9").Value = 208000
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Jan 2026"
$ws.Range("G9").Value = 364484.6743295019
$ws.Range("H9").Value = -1000
$ws.Range("I9").Value = -0.004784688995215311
